$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 and B1 hold text-typed numbers (quoted inline strings in the sheet XML).
# Setting Range.Value to a numeric-looking string auto-converts it to a
# Number cell, so force text via NumberFormat "@" and then restore the
# original cell format (copied from the untouched C1, which keeps style
# index 1) so the style id does not drift.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "2"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2.1"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# C1 stays "0" - untouched.

# Row 2 numeric values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1

# Drop column D (header + data rows) and row 3 (old data no longer present)
$ws.Range("D1:D3").Delete() | Out-Null
$ws.Range("A3:C3").Delete() | Out-Null
